$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ----------------------------------------------------------------------
# New TA meeting attendance block for 12/1/2020 (rows 25-31), mirroring
# the layout already used for the other meeting blocks above it.
# ----------------------------------------------------------------------

# Row 25: date header. B25 holds the meeting date; C25:N25 are blank,
# center-aligned cells (same look as the "J17:M17" filler cells).
$ws.Range("B25:N25").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B25").NumberFormat = "d-mmm"
$ws.Range("B25").Value = 44166

# Row 26: meeting type / TA marker.
$ws.Range("B26:N26").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A26").Value = "Meeting Type"
$ws.Range("B26").Value = "TA"

# Rows 27-31: team member attendance cells, colored like the legend
# (green == present). Row 27 (Yao Siyu) mirrors the "planned absence"
# light-green highlight used elsewhere (e.g. D19); rows 28-31 mirror the
# default "present" green fill (e.g. C20:C23).
$ws.Range("A27").Value = "Yao Siyu"
$ws.Range("B27").Interior.Color = 5296274   # RGB(146,208,80) light green

$ws.Range("A28").Value = "Rivas, Madison"
$ws.Range("B28").Interior.Color = 5287936   # RGB(0,176,80) green

$ws.Range("A29").Value = "Alvaro Santillan"
$ws.Range("B29").Interior.Color = 5287936   # RGB(0,176,80) green

$ws.Range("A30").Value = "Liyang Ru"
$ws.Range("B30").Interior.Color = 5287936   # RGB(0,176,80) green

$ws.Range("A31").Value = "Guangshi Xu"
$ws.Range("B31").Interior.Color = 5287936   # RGB(0,176,80) green

# Scroll the view down to the new block and leave the selection where the
# author left off (matches the sheetView's topLeftCell="A8" / activeCell
# "C30" saved in the workbook).
$ws.Range("C30").Select()
$excel.ActiveWindow.ScrollRow = 8
